# "funciona bien las rutas de accion"
# The "Provincia" column (D) for postal codes 1000-1599 (rows 2-459) was
# mis-tagged as "CAPITAL FEDERAL" — correct it to "BUENOS AIRES", matching
# the rest of the sheet (rows 460+ already read "BUENOS AIRES").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bulk-update column D (Provincia) for rows 2 through 459.
$ws.Range("D2:D459").Value = "BUENOS AIRES"

# Restore the author's on-save view state: scrolled down to row 2132,
# with A2143 as the active cell/selection.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 2132
$win.ScrollColumn = 1
$ws.Range("A2143").Select() | Out-Null
